$d = $word.ActiveDocument

# Locate the paragraph holding the user-doc opening marker, currently split
# across two runs: "{m" and ":userdoc 'zone1'}". We need it re-split into
# four runs: "{", "m", ":userdoc 'zone1'", "}" (matching the
# TokenIteratorFieldRewriterSplit output), preserving the paragraph's
# existing formatting/rsid attributes.

$p = $null
foreach ($candidate in $d.Paragraphs) {
    if ($candidate.Range.Text.StartsWith("{m:userdoc")) {
        $p = $candidate
        break
    }
}
if ($p -eq $null) {
    $p = $d.Paragraphs.Item(2)
}
$r = $p.Range

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParaXml = '<w:p ' + $w + ' w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979">' +
    '<w:r w:rsidR="00DE6D5A"><w:t>{</w:t></w:r>' +
    '<w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r>' +
    '<w:r w:rsidR="00833091"><w:t>:userdoc ' + "'" + 'zone1' + "'" + '</w:t></w:r>' +
    '<w:r w:rsidR="00833091"><w:t xml:space="preserve">}</w:t></w:r>' +
    '</w:p>'

$null = $r.InsertXML($newParaXml)
